# Apply the edit described in the diff:
# - Rows 7-11 are removed (dimension shrinks from A1:B11 to A1:B6)
# - Remaining rows 2-6 get updated A/B values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 7 through 11 entirely (shifts nothing below, they are the last rows)
$ws.Range("A7:B11").EntireRow.Delete() | Out-Null

# Update the surviving data rows with their new values
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = 25

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 23

$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 17

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 17

$ws.Cells.Item(6, 1).Value = 0
$ws.Cells.Item(6, 2).Value = 17
